# Apply the "refactored" location data for the Sabra and Shatila massacre sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set (header stays the same, data rows replaced/truncated from 16 to 10 rows)
$data = @(
    @("Name", "Subdivision type", "Subdivision name"),
    @("Beirut", "Country", "Lebanon"),
    @("Iraq", "-", "-"),
    @("Lebanon", "-", "-"),
    @("Nahariyya", "Country", "ISR"),
    @("Sabra", "Country", "Lebanon"),
    @("Saliha", "Geopolitical entity", "Mandatory Palestine"),
    @("Singapore", "-", "-"),
    @("Israeli occupation of Southern Lebanon", "-", "-"),
    @("Southern Lebanon", "-", "-"),
    @("Beirut", "Country", "Lebanon")
)

# Clear out the old used range contents first so leftover rows (12-17) are
# removed, while keeping existing formatting (e.g. the bold header row) intact.
$usedRange = $ws.UsedRange
$usedRange.ClearContents()

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
